# This script appends a new Q&A block (question 12: accept filename and
# print its extension) after the existing "4. Store a list..." block
# that ends with "The Number of occurrence of a: 2".
#
# Strategy: locate the paragraph that ends with the text
# "The Number of occurrence of a: 2", then splice in 8 new paragraphs
# after it (matching the pPr/rPr formatting already used throughout the
# document), writing each paragraph's full OOXML (runs + proofErr markers)
# via Range.InsertXML so the result matches Word's own "type it in" output.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Common run properties used by every run in this block.
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="28"/></w:rPr>'

# pPr variants used across the new paragraphs.
$pPrWithBefore = '<w:pPr><w:spacing w:before="184" w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="140"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="28"/></w:rPr></w:pPr>'
$pPrNoBefore   = '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="140"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="28"/></w:rPr></w:pPr>'

function New-Run($text, [switch]$preserve) {
    if ($preserve) {
        return '<w:r>' + $rPr + '<w:t xml:space="preserve">' + $text + '</w:t></w:r>'
    } else {
        return '<w:r>' + $rPr + '<w:t>' + $text + '</w:t></w:r>'
    }
}

# Locate the anchor paragraph: "The Number of occurrence of a: 2"
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd() -eq "The Number of occurrence of a: 2") {
        $anchor = $cand
    }
}
if ($anchor -eq $null) {
    throw "Could not find anchor paragraph 'The Number of occurrence of a: 2'"
}

$anchorIndex = $anchor.Index

# Create 8 blank paragraphs right after the anchor paragraph; each
# InsertParagraphAfter call adds one new empty paragraph immediately
# following the (still-collapsed) range, so repeating it keeps stacking
# them in order right after the anchor.
$rng = $anchor.Range
$rng.Collapse(0)
for ($n = 0; $n -lt 8; $n++) {
    $rng.InsertParagraphAfter()
}

# Paragraph 1: blank spacer paragraph (before=184)
$p1 = $d.Paragraphs.Item($anchorIndex + 1)
$p1.Range.InsertXML('<w:p ' + $wNs + '>' + $pPrWithBefore + '</w:p>')

# Paragraph 2: "12.Accept a file name from user and print extension of that"
$p2 = $d.Paragraphs.Item($anchorIndex + 2)
$xml2 = '<w:p ' + $wNs + '>' + $pPrNoBefore `
    + (New-Run "12") `
    + '<w:proofErr w:type="gramStart"/>' `
    + (New-Run ".Accept") `
    + '<w:proofErr w:type="gramEnd"/>' `
    + (New-Run " a file name from user and print extension of that" -preserve) `
    + '</w:p>'
$p2.Range.InsertXML($xml2)

# Paragraph 3: str1=input("Enter the filename:")
$p3 = $d.Paragraphs.Item($anchorIndex + 3)
$xml3 = '<w:p ' + $wNs + '>' + $pPrWithBefore `
    + (New-Run "str1=") `
    + '<w:proofErr w:type="gramStart"/>' `
    + (New-Run "input(") `
    + '<w:proofErr w:type="gramEnd"/>' `
    + (New-Run '"Enter the filename:")') `
    + '</w:p>'
$p3.Range.InsertXML($xml3)

# Paragraph 4: str2=str1.split(".")
$p4 = $d.Paragraphs.Item($anchorIndex + 4)
$xml4 = '<w:p ' + $wNs + '>' + $pPrWithBefore `
    + (New-Run "str2=") `
    + '<w:proofErr w:type="gramStart"/>' `
    + (New-Run "str1.split(") `
    + '<w:proofErr w:type="gramEnd"/>' `
    + (New-Run '".")') `
    + '</w:p>'
$p4.Range.InsertXML($xml4)

# Paragraph 5: print("The File Extension:",str2[1])
$p5 = $d.Paragraphs.Item($anchorIndex + 5)
$xml5 = '<w:p ' + $wNs + '>' + $pPrWithBefore `
    + '<w:proofErr w:type="gramStart"/>' `
    + (New-Run "print(") `
    + '<w:proofErr w:type="gramEnd"/>' `
    + (New-Run '"The File Extension:",str2[1])') `
    + '</w:p>'
$p5.Range.InsertXML($xml5)

# Paragraph 6: output:
$p6 = $d.Paragraphs.Item($anchorIndex + 6)
$runOutput = '<w:r>' + $rPr + '<w:lastRenderedPageBreak/><w:t>output</w:t></w:r>'
$xml6 = '<w:p ' + $wNs + '>' + $pPrWithBefore `
    + '<w:proofErr w:type="gramStart"/>' `
    + $runOutput `
    + '<w:proofErr w:type="gramEnd"/>' `
    + (New-Run ":") `
    + '</w:p>'
$p6.Range.InsertXML($xml6)

# Paragraph 7: Enter the filename:hiii.java
$p7 = $d.Paragraphs.Item($anchorIndex + 7)
$xml7 = '<w:p ' + $wNs + '>' + $pPrWithBefore `
    + (New-Run "Enter the " -preserve) `
    + '<w:proofErr w:type="spellStart"/>' `
    + (New-Run "filename") `
    + '<w:proofErr w:type="gramStart"/>' `
    + (New-Run ":hiii.java") `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:proofErr w:type="gramEnd"/>' `
    + '</w:p>'
$p7.Range.InsertXML($xml7)

# Paragraph 8: The File Extension: java  (bookmark stays attached after this
# run, matching where Word leaves the _GoBack mark after the last edit)
$p8 = $d.Paragraphs.Item($anchorIndex + 8)
$xml8 = '<w:p ' + $wNs + '>' + $pPrWithBefore `
    + (New-Run "The File Extension: java") `
    + '<w:bookmarkStart w:id="100" w:name="_GoBack"/><w:bookmarkEnd w:id="100"/>' `
    + '</w:p>'
$p8.Range.InsertXML($xml8)

# Remove the old bookmark that was left at the end of the anchor paragraph,
# since Word only keeps a single _GoBack bookmark (at the most recent edit).
$bm = $d.Bookmarks
for ($i = 1; $i -le $bm.Count; $i++) {
    $item = $bm.Item($i)
    if ($item.Name -eq "_GoBack" -and $item.Range.Start -lt $p1.Range.Start) {
        $item.Delete()
        break
    }
}
